{"js": "// Update the three-digit x one-digit multiplication prompts in the practice\n// table to a new generated set of problems (old => new).\nconst replacements = [\n  [\"153\u00d77=\", \"120\u00d75=\"],\n  [\"506\u00d79=\", \"345\u00d79=\"],\n  [\"993\u00d77=\", \"316\u00d73=\"],\n  [\"124\u00d76=\", \"219\u00d78=\"],\n  [\"220\u00d74=\", \"654\u00d77=\"],\n  [\"517\u00d78=\", \"564\u00d79=\"],\n  [\"235\u00d74=\", \"491\u00d79=\"],\n  [\"239\u00d79=\", \"207\u00d75=\"],\n  [\"668\u00d74=\", \"525\u00d73=\"],\n  [\"758\u00d77=\", \"115\u00d75=\"],\n  [\"800\u00d74=\", \"896\u00d75=\"],\n  [\"141\u00d77=\", \"839\u00d79=\"],\n  [\"971\u00d72=\", \"184\u00d74=\"],\n  [\"763\u00d79=\", \"442\u00d79=\"],\n  [\"741\u00d75=\", \"332\u00d74=\"],\n  [\"364\u00d78=\", \"249\u00d73=\"],\n  [\"852\u00d76=\", \"361\u00d75=\"],\n  [\"253\u00d77=\", \"692\u00d75=\"],\n  [\"991\u00d75=\", \"911\u00d78=\"],\n  [\"532\u00d79=\", \"801\u00d74=\"],\n  [\"193\u00d73=\", \"782\u00d72=\"],\n  [\"493\u00d77=\", \"165\u00d73=\"],\n  [\"808\u00d75=\", \"351\u00d79=\"],\n  [\"610\u00d78=\", \"264\u00d76=\"],\n  [\"640\u00d73=\", \"101\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit x one-digit multiplication prompts in the practice\n# table to a new generated set of problems (old => new):\n#   153\u00d77= -> 120\u00d75=   506\u00d79= -> 345\u00d79=   993\u00d77= -> 316\u00d73=\n#   124\u00d76= -> 219\u00d78=   220\u00d74= -> 654\u00d77=   517\u00d78= -> 564\u00d79=\n#   235\u00d74= -> 491\u00d79=   239\u00d79= -> 207\u00d75=   668\u00d74= -> 525\u00d73=\n#   758\u00d77= -> 115\u00d75=   800\u00d74= -> 896\u00d75=   141\u00d77= -> 839\u00d79=\n#   971\u00d72= -> 184\u00d74=   763\u00d79= -> 442\u00d79=   741\u00d75= -> 332\u00d74=\n#   364\u00d78= -> 249\u00d73=   852\u00d76= -> 361\u00d75=   253\u00d77= -> 692\u00d75=\n#   991\u00d75= -> 911\u00d78=   532\u00d79= -> 801\u00d74=   193\u00d73= -> 782\u00d72=\n#   493\u00d77= -> 165\u00d73=   808\u00d75= -> 351\u00d79=   610\u00d78= -> 264\u00d76=\n#   640\u00d73= -> 101\u00d76=\n\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n    \"153\u00d77=\" = \"120\u00d75=\"\n    \"506\u00d79=\" = \"345\u00d79=\"\n    \"993\u00d77=\" = \"316\u00d73=\"\n    \"124\u00d76=\" = \"219\u00d78=\"\n    \"220\u00d74=\" = \"654\u00d77=\"\n    \"517\u00d78=\" = \"564\u00d79=\"\n    \"235\u00d74=\" = \"491\u00d79=\"\n    \"239\u00d79=\" = \"207\u00d75=\"\n    \"668\u00d74=\" = \"525\u00d73=\"\n    \"758\u00d77=\" = \"115\u00d75=\"\n    \"800\u00d74=\" = \"896\u00d75=\"\n    \"141\u00d77=\" = \"839\u00d79=\"\n    \"971\u00d72=\" = \"184\u00d74=\"\n    \"763\u00d79=\" = \"442\u00d79=\"\n    \"741\u00d75=\" = \"332\u00d74=\"\n    \"364\u00d78=\" = \"249\u00d73=\"\n    \"852\u00d76=\" = \"361\u00d75=\"\n    \"253\u00d77=\" = \"692\u00d75=\"\n    \"991\u00d75=\" = \"911\u00d78=\"\n    \"532\u00d79=\" = \"801\u00d74=\"\n    \"193\u00d73=\" = \"782\u00d72=\"\n    \"493\u00d77=\" = \"165\u00d73=\"\n    \"808\u00d75=\" = \"351\u00d79=\"\n    \"610\u00d78=\" = \"264\u00d76=\"\n    \"640\u00d73=\" = \"101\u00d76=\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
